$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add 10 new AI Engineer job listings (rows 87-96)
# Row 87: Palantir - Forward Deployed AI Engineer
$ws.Range('A87').Value = 'Palantir'
$ws.Range('B87').Value = 'Forward Deployed AI Engineer'
$ws.Range('C87').Value = 'Big Tech'
$ws.Range('D87').Value = 'New York, NY / Washington, DC'
$ws.Range('E87').Value = 'AI Engineer'
$ws.Range('F87').Value = 'https://jobs.lever.co/palantir/636fc05c-d348-4a06-be51-597cb9e07488'
$ws.Range('G87').Value = 'Strong engineering background in CS, Math, Software Engineering, Physics, or ML. Deep understanding of the Gen AI landscape.'
$ws.Range('H87').Value = 'Build end-to-end LLM workflows at scale for enterprise customers. Own Gen AI strategy and implementation. Work directly with customers to solve real-world problems — role resembles a hands-on AI startup CTO. Comfort working in dynamic environments with evolving objectives.'
$ws.Range('I87').Value = 'Python; LLMs, prompt engineering, agent development; ML fundamentals (evaluation, training, problem decomposition); production Gen AI systems'

# Row 88: Quora (Poe) - AI Engineer New Grad 2025-2026
$ws.Range('A88').Value = 'Quora (Poe)'
$ws.Range('B88').Value = 'AI Engineer New Grad 2025-2026'
$ws.Range('C88').Value = 'Unicorn'
$ws.Range('D88').Value = 'Remote'
$ws.Range('E88').Value = 'AI Engineer'
$ws.Range('F88').Value = 'https://jobs.ashbyhq.com/quora/6df58d3e-855a-423e-99fd-a56ac8824b34'
$ws.Range('G88').Value = '2025 or Summer 2026 graduate with B.S., M.S., or Ph.D. in CS, Engineering, or related technical field.'
$ws.Range('H88').Value = 'Work on prompt engineering, retrieval-augmented generation, and agentic workflow optimization. Improve existing applied AI systems and identify new opportunities to apply emerging AI capabilities. Take end-to-end ownership from prototyping, data pipelines, model optimization/evaluation to reliable deployment at scale.'
$ws.Range('I88').Value = 'Python, TypeScript; LLM prompt engineering; RAG; agentic workflows; model evaluation and optimization; data pipelines'

# Row 89: Databricks - AI Engineer - FDE (Forward Deployed Engineer)
$ws.Range('A89').Value = 'Databricks'
$ws.Range('B89').Value = 'AI Engineer - FDE (Forward Deployed Engineer)'
$ws.Range('C89').Value = 'Unicorn'
$ws.Range('D89').Value = 'Remote / US'
$ws.Range('E89').Value = 'AI Engineer'
$ws.Range('F89').Value = 'https://www.databricks.com/company/careers/professional-services-operations/ai-engineer---fde-forward-deployed-engineer-8024010002'
$ws.Range('G89').Value = 'Graduate degree in CS, Engineering, Statistics, Operations Research, or equivalent practical experience. Passion for driving business value through AI.'
$ws.Range('H89').Value = 'Develop cutting-edge GenAI solutions for customers using latest techniques from Mosaic AI Research. Embed with customer teams from technical ICs to executives. Contribute accelerators, frameworks, and best practices that scale across accounts and influence product roadmap. Travel up to 50%.'
$ws.Range('I89').Value = 'Python, SQL, Java/Scala, JavaScript/TypeScript; AWS/Azure/GCP; Apache Spark; Databricks Intelligence Platform; OpenAI/Anthropic/Gemini APIs; production ML deployments; distributed datasets'

# Row 90: Scale AI - Applied AI Engineer, Enterprise GenAI
$ws.Range('A90').Value = 'Scale AI'
$ws.Range('B90').Value = 'Applied AI Engineer, Enterprise GenAI'
$ws.Range('C90').Value = 'Unicorn'
$ws.Range('D90').Value = 'San Francisco, CA / New York, NY'
$ws.Range('E90').Value = 'AI Engineer'
$ws.Range('F90').Value = 'https://scale.com/careers/4514173005'
$ws.Range('G90').Value = 'Bachelor''s in CS, Mathematics, or related quantitative field. Strong Python proficiency. Cloud platform experience (AWS or GCP).'
$ws.Range('H90').Value = 'Own, plan, and optimize AI behind enterprise customers'' deepest technical problems. Build advanced AI agents with multimodal and tool-calling capabilities on Scale''s Generative Platform. Convert business requirements into technical implementations. Write and debug production code across company and customer environments.'
$ws.Range('I90').Value = 'Python, NumPy, Pandas; AWS, GCP; LLMs, generative AI applications; AI agents with tool-calling; production ML model development'

# Row 91: Scale AI - Forward Deployed AI Engineer, Enterprise
$ws.Range('A91').Value = 'Scale AI'
$ws.Range('B91').Value = 'Forward Deployed AI Engineer, Enterprise'
$ws.Range('C91').Value = 'Unicorn'
$ws.Range('D91').Value = 'San Francisco, CA / New York, NY'
$ws.Range('E91').Value = 'AI Engineer'
$ws.Range('F91').Value = 'https://scale.com/careers/4597399005'
$ws.Range('G91').Value = '4+ years software engineering experience. Production Python expertise with LangChain, LlamaIndex, HuggingFace, OpenAI API. Cloud platform experience (AWS, GCP, Azure).'
$ws.Range('H91').Value = 'Partner with enterprise clients on infrastructure and data pipeline requirements. Develop production-grade agents for customer support, analysis, content generation, automation. Architect multi-agent systems across models and data sources. Implement RAG systems, fine-tuning pipelines, and human-in-the-loop feedback. Build data connectors and ETL pipelines. Serve as primary technical contact for strategic accounts.'
$ws.Range('I91').Value = 'Python; LangChain, LlamaIndex, HuggingFace, OpenAI API; RAG, embeddings, vector databases; Docker, Kubernetes, CI/CD; Terraform/IaC; multi-agent systems; prompt engineering'

# Row 92: Anthropic - Forward Deployed Engineer, Applied AI
$ws.Range('A92').Value = 'Anthropic'
$ws.Range('B92').Value = 'Forward Deployed Engineer, Applied AI'
$ws.Range('C92').Value = 'Unicorn'
$ws.Range('D92').Value = 'San Francisco, CA'
$ws.Range('E92').Value = 'AI Engineer'
$ws.Range('F92').Value = 'https://job-boards.greenhouse.io/anthropic/jobs/4985877008'
$ws.Range('G92').Value = '4+ years experience in a technical, customer-facing role (Forward Deployed Engineer, Solutions Engineer) or as a Software Engineer with consulting experience.'
$ws.Range('H92').Value = 'Production experience with LLMs including advanced prompt engineering and agent development. Work with customers to implement AI solutions that solve their real business problems. Bridge gap between Anthropic''s AI capabilities and enterprise customer needs.'
$ws.Range('I92').Value = 'Python; TypeScript, Java; LLMs, advanced prompt engineering; agent development; customer-facing technical implementation'

# Row 93: Snowflake - Applied AI Engineer
$ws.Range('A93').Value = 'Snowflake'
$ws.Range('B93').Value = 'Applied AI Engineer'
$ws.Range('C93').Value = 'Big Tech'
$ws.Range('D93').Value = 'San Mateo, CA / Remote'
$ws.Range('E93').Value = 'AI Engineer'
$ws.Range('F93').Value = 'https://careers.snowflake.com/us/en/ai-ml-engineering'
$ws.Range('G93').Value = 'Bachelor''s in CS, Engineering, or related field. 2+ years professional software engineering experience. Experience in a customer-facing role (solutions architect, sales engineer, or professional services).'
$ws.Range('H93').Value = 'Hands-on builder and critical technical partner to strategic customers at the forefront of enterprise AI. Build, evaluate, and tune applications and pipelines involving ML models or data-intensive systems. Tackle complex and ambiguous technical challenges leveraging cutting-edge research and AI.'
$ws.Range('I93').Value = 'Python; Snowpark, pandas, NumPy; ML model evaluation and tuning; data-intensive pipelines; customer-facing technical delivery'

# Row 94: Amazon - Software Dev Engineer Intern - AI/ML (Summer 2026), AGI
$ws.Range('A94').Value = 'Amazon'
$ws.Range('B94').Value = 'Software Dev Engineer Intern - AI/ML (Summer 2026), AGI'
$ws.Range('C94').Value = 'FAANG'
$ws.Range('D94').Value = 'Seattle, WA / Sunnyvale, CA'
$ws.Range('E94').Value = 'AI Engineer'
$ws.Range('F94').Value = 'https://www.amazon.jobs/en/jobs/3121382/software-dev-engineer-intern-ai-ml-summer-2026-artificial-general-intelligence-agi'
$ws.Range('G94').Value = 'Enrolled in B.S. or above in CS, Computer Engineering, Data Science, Electrical Engineering, or related STEM field. 18 years or older.'
$ws.Range('H94').Value = 'Experience programming with at least one modern language (Java, C++, Python). Experience with AI/ML technologies. Work on Amazon''s Artificial General Intelligence (AGI) organization building next-generation AI products and services.'
$ws.Range('I94').Value = 'Python, Java, C++; AI/ML technologies; cloud platforms (AWS); debugging complex systems; software development lifecycle'

# Row 95: Google - AI Developer Engineer, Cloud AI
$ws.Range('A95').Value = 'Google'
$ws.Range('B95').Value = 'AI Developer Engineer, Cloud AI'
$ws.Range('C95').Value = 'FAANG'
$ws.Range('D95').Value = 'Sunnyvale, CA / New York, NY / Kirkland, WA'
$ws.Range('E95').Value = 'AI Engineer'
$ws.Range('F95').Value = 'https://www.google.com/about/careers/applications/jobs/results/102301475132056262-ai-developer-engineer/'
$ws.Range('G95').Value = 'Bachelor''s degree or equivalent practical experience. 2 years experience with software development in Java, Python, or C++, or 1 year with an advanced degree.'
$ws.Range('H95').Value = '1 year experience with ML infrastructure (model deployment, evaluation, optimization, data processing, debugging). Experience with GenAI techniques (LLMs, multi-modal, large vision models) or GenAI-related concepts (language modeling, computer vision). Build and deploy AI solutions on Google Cloud for enterprise customers.'
$ws.Range('I95').Value = 'Python, Java, C++; TensorFlow, PyTorch; ML infrastructure; LLMs, multi-modal models; Google Cloud Platform; model deployment and evaluation'

# Row 96: ServiceNow - Machine Learning Engineer Intern - Summer 2026
$ws.Range('A96').Value = 'ServiceNow'
$ws.Range('B96').Value = 'Machine Learning Engineer Intern - Summer 2026'
$ws.Range('C96').Value = 'Big Tech'
$ws.Range('D96').Value = 'Santa Clara, CA (Hybrid)'
$ws.Range('E96').Value = 'AI Engineer'
$ws.Range('F96').Value = 'https://careers.servicenow.com/jobs/744000080670612/machine-learning-engineer-intern-summer-2026/'
$ws.Range('G96').Value = 'Current enrollment in full-time Bachelor''s or Master''s in AI, ML, Data Science, CS, or related discipline. Strong theoretical and practical knowledge of ML/deep learning.'
$ws.Range('H96').Value = 'Work alongside engineers to enhance ServiceNow''s NLP capabilities using advanced deep learning algorithms. Apply ML to automate enterprise workflows at scale. 12-week internship (May-August or June-September).'
$ws.Range('I96').Value = 'Java, Python, JavaScript, SQL; NumPy, Pandas, Scikit-learn, TensorFlow/PyTorch; supervised/unsupervised learning; model evaluation and optimization; NLP/deep learning'

